# Add "NA" values in column E (duplicate_image_filename) for rows 2-21,
# matching the header already present in E1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
